# Logged Week 15 and simulated Week 16
# Adds a new player column "J.Sprinkle" (with placeholder value "n") to the
# Yards Data table on both the "Rushing" and "Receiving" sheets, mirroring
# the existing header formatting used for the other player columns.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $header = $ws.Range("U1")

    # Match the bordered / bold / centered header style used by the other
    # player-name header cells (e.g. T1) before writing the new header text.
    $header.Font.Bold = $true
    $header.HorizontalAlignment = -4108   # xlCenter
    $header.VerticalAlignment = -4160     # xlTop
    $header.Borders.LineStyle = 1         # xlContinuous
    $header.Value = "J.Sprinkle"

    # New player's data row placeholder value, same as the rest of row 2.
    $ws.Range("U2").Value = "n"
}
